# "update on ahp algo: 1. file transfer, 2. eigen value"
#
# The old "Level1" sheet (long pairwise-comparison labels + a stray marker
# cell in C10 + custom column widths) is replaced with a freshly rebuilt
# "Level1" sheet that uses the shorter label names, keeps only the real
# data (A1:F2), and sits right after "Level2_Econ" in the tab order. The
# previously-active tab ("Level2_Technical") hands off "active" status to
# the rebuilt "Level1" sheet.

$wb = $excel.ActiveWorkbook

# Drop the old Level1 sheet (long labels, stray C10 "`" cell, custom
# column widths) entirely...
$old = $wb.Worksheets.Item("Level1")
$old.Delete() | Out-Null

# ...and recreate it fresh. Worksheets.Add() (no args) inserts right
# before the currently active sheet ("Level2_Technical"), i.e. right
# after "Level2_Econ" - exactly the new tab position we want.
$newLevel1 = $wb.Worksheets.Add()
$newLevel1.Name = "Level1"

# Make sure it's positioned immediately before Level2_Technical (in case
# sheet order wasn't already right).
$newLevel1.Move($wb.Worksheets.Item("Level2_Technical"))

# Re-fetch a live reference after the Move (the old COM reference can
# otherwise point at stale workbook state).
$ws = $wb.Worksheets.Item("Level1")

$ws.Range("A1").Value = "Tech_Econ"
$ws.Range("B1").Value = "Tech_Infra"
$ws.Range("C1").Value = "Tech_serviceStandard"
$ws.Range("D1").Value = "Econ_Infra"
$ws.Range("E1").Value = "Econ_serviceStandard"
$ws.Range("F1").Value = "Infra_serviceStandard"

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = -3
$ws.Range("F2").Value = -3

$ws.Range("A2").Select()

# This sheet becomes the active tab (matches activeTab moving to index 1).
$ws.Activate()
